$d = $word.ActiveDocument

function Replace-InParagraph($index, $findText, $replaceText) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

# 1. Para 16: "Actividad de reflexión sobre el género narrativo"
#    -> "Actividad para afianzar el conocimiento sobre el género narrativo"
Replace-InParagraph 16 "de reflexión" "para afianzar el conocimiento"

# 2. Para 99: "2-Medio" -> "2"
Replace-InParagraph 99 "2-Medio" "2"

# 3. Para 110: "S" -> "“S”"
Replace-InParagraph 110 "S" "“S”"

# 4. Para 114: reword the "Realiza la siguiente actividad..." instructions
Replace-InParagraph 114 "Si es necesario puedes entregar tu respuesta a mano durante la clase, o por email a tu profesor para que pueda evaluarla." "Si es necesario, puedes entregarle tu respuesta a tu docente, a mano o por email, para que pueda evaluarla."

# 5. Para 132: "Utiliza las palabras ..." quoted words -> unquoted + wording change
Replace-InParagraph 132 "Utiliza las palabras “sierpe”, “almohada”, “carta”, “tornamesa”, “alcaraván”; para realizar un cuento breve." "Utiliza las palabras sierpe, almohada, carta, tornamesa y alcaraván dentro de un cuento breve de tu inspiración."

# 6. Para 136: "2- Medio" -> "2"
Replace-InParagraph 136 "2- Medio" "2"

# 7. Para 153: "... y contesta qué elementos del género narrativo se encuentran presentes en él."
Replace-InParagraph 153 "y contesta qué elementos del género narrativo se encuentran presentes en él." "e identifica los elementos del género narrativo que se encuentran presentes."

# 8. Para 157: "1-Fácil" -> "1"
Replace-InParagraph 157 "1-Fácil" "1"

# 9. Para 170: leading "-" -> en dash "–" (only first char of the paragraph)
Replace-InParagraph 170 "-Así estaban las cosas." "–Así estaban las cosas."

# 10. Para 171: "-Bien -le dije-, tengo que subir a enviar unos cables."
Replace-InParagraph 171 "-Bien -le dije-, tengo que subir a enviar unos cables." "–Bien –le dije–, tengo que subir a enviar unos cables."

# 11. Para 172: "-¿De veras?"
Replace-InParagraph 172 "-¿De veras?" "–¿De veras?"

# 12. Para 173: "-Sí, tengo que mandar unos cables."
Replace-InParagraph 173 "-Sí, tengo que mandar unos cables." "–Sí, tengo que mandar unos cables."

# 13. Para 174: "-¿Te importa si subo y me quedo por ahí en tu despacho?"
Replace-InParagraph 174 "-¿Te importa si subo y me quedo por ahí en tu despacho?" "–¿Te importa si subo y me quedo por ahí en tu despacho?"

# 14. Para 175: "-No, sube."
Replace-InParagraph 175 "-No, sube." "–No, sube."

# 15. Para 177: "Ernest Hemingway. Fiesta." -> "Ernest Hemingway, Fiesta "
Replace-InParagraph 177 "Ernest Hemingway. Fiesta." "Ernest Hemingway, Fiesta "
